$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENEFITS")

# A3: replace the "Benefits for the period ..." formula with a static,
# already-computed text value for this run (25/09/2025 to 01/10/2025).
$ws.Range("A3").Value = "Benefits period: 25/09/2025 to 01/10/2025"

# D4:D19 - wire each row up to the per-client external workbook reference
# (these used to be blank placeholders / one-off hard-coded numbers; they
# now pull live from the named tabs in the linked Client Funds workbook).
$ws.Range("D4").Formula  = "=CLAGUE!D643"
$ws.Range("D5").Formula  = "=COLLISTER!D147"
$ws.Range("D6").Formula  = "=CORKILL!D561"
$ws.Range("D7").Formula  = "=DURRANT!D105"
$ws.Range("D8").Formula  = "=DYER!D587"
$ws.Range("D9").Formula  = "=JACKSON!D623"
$ws.Range("D10").Formula = "=JONES!D7"
$ws.Range("D11").Formula = "=LEWIS!D610"
$ws.Range("D12").Formula = "=MCLAREN!D468"
$ws.Range("D13").Formula = "=PATTON!D565"
$ws.Range("D14").Formula = "=PERRY!D606"
$ws.Range("D15").Formula = "=SAYLE!D321"
$ws.Range("D16").Formula = "=SHIMMIN!D87"
$ws.Range("D17").Formula = "=SMITH!D598"
$ws.Range("D18").Formula = "=WARD!D607"

# D19 (TOTAL BALANCE IN ACCOUNTS) now pulls from WEST instead of summing.
$ws.Range("D19").Formula = "=WEST!D66"

# Editing the formula cells above nudges Excel's auto row-height for a few
# rows; pin them back to their original explicit heights.
$ws.Rows.Item(4).RowHeight = 6.75
$ws.Rows.Item(14).RowHeight = 14.1
$ws.Rows.Item(18).RowHeight = 14.1

# Materialise row 24 (previously skipped between rows 23 and 25).
$ws.Rows.Item(24).RowHeight = $ws.Rows.Item(24).RowHeight
